$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 2850.5
$ws_ALC.Range("I4").Value = 2701
$ws_ALC.Range("J4").Value = 3000
$ws_ALC.Range("K4").Value = 2701
$ws_ALC.Range("L4").Value = 3000
$ws_ALC.Range("M4").Value = -2587
$ws_ALC.Range("N4").Value = -3228

$ws_ALC.Range("H9").Value = 149.57143
$ws_ALC.Range("I9").Value = 169.4
$ws_ALC.Range("J9").Value = 100
$ws_ALC.Range("K9").Value = 169.4
$ws_ALC.Range("L9").Value = 100
$ws_ALC.Range("M9").Value = -0.4000000000000057
$ws_ALC.Range("N9").Value = -438

$ws_ALC.Range("H62").Value = 4981.8335
$ws_ALC.Range("I62").Value = 5918.6
$ws_ALC.Range("J62").Value = 298
$ws_ALC.Range("K62").Value = 5918.6
$ws_ALC.Range("L62").Value = 298
$ws_ALC.Range("M62").Value = -5294.6
$ws_ALC.Range("N62").Value = -1546

$ws_ALC.Range("H65").Value = 4981.8335
$ws_ALC.Range("I65").Value = 5918.6
$ws_ALC.Range("J65").Value = 298
$ws_ALC.Range("K65").Value = 29593
$ws_ALC.Range("L65").Value = 1490
$ws_ALC.Range("M65").Value = -26473
$ws_ALC.Range("N65").Value = -7730

$ws_ALC.Range("H113").Value = 27263.75

$ws_ALC.Range("H137").Value = 46736.09
$ws_ALC.Range("J137").Value = 92503.37
$ws_ALC.Range("L137").Value = 277510.11
$ws_ALC.Range("N137").Value = -282610.11

$ws_ALC.Range("H138").Value = 1743.6262
$ws_ALC.Range("I138").Value = 918.4103
$ws_ALC.Range("J138").Value = 2280.0166
$ws_ALC.Range("K138").Value = 2755.2309
$ws_ALC.Range("L138").Value = 6840.0498
$ws_ALC.Range("M138").Value = 2384.7691
$ws_ALC.Range("N138").Value = -17120.0498

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4898.089
$ws_ARM.Range("I32").Value = 2909.338
$ws_ARM.Range("J32").Value = 14096.0625
$ws_ARM.Range("K32").Value = 2909.338
$ws_ARM.Range("L32").Value = 14096.0625
$ws_ARM.Range("M32").Value = -2622.338
$ws_ARM.Range("N32").Value = -14670.0625

$ws_ARM.Range("H45").Value = 1441.3158
$ws_ARM.Range("I45").Value = 1235.6364
$ws_ARM.Range("K45").Value = 1235.6364
$ws_ARM.Range("M45").Value = -858.6364000000001

$ws_ARM.Range("H122").Value = 2467.889
$ws_ARM.Range("I122").Value = 1526.375
$ws_ARM.Range("J122").Value = 10000
$ws_ARM.Range("K122").Value = 4579.125
$ws_ARM.Range("L122").Value = 30000
$ws_ARM.Range("M122").Value = -2129.125
$ws_ARM.Range("N122").Value = -34900

$ws_ARM.Range("H132").Value = 1482.3857
$ws_ARM.Range("I132").Value = 1206.1459
$ws_ARM.Range("K132").Value = 3618.4377
$ws_ARM.Range("M132").Value = -1088.4377

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 2960.5
$ws_BSM.Range("I20").Value = 3123.6365
$ws_BSM.Range("J20").Value = 2704.1428
$ws_BSM.Range("K20").Value = 3123.6365
$ws_BSM.Range("L20").Value = 2704.1428
$ws_BSM.Range("M20").Value = -2876.6365
$ws_BSM.Range("N20").Value = -3198.1428

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H107").Value = 784.03125
$ws_CRP.Range("I107").Value = 683.5417
$ws_CRP.Range("K107").Value = 683.5417
$ws_CRP.Range("M107").Value = 1236.4583

$ws_CRP.Range("H122").Value = 4122.143
$ws_CRP.Range("I122").Value = 2124.875
$ws_CRP.Range("J122").Value = 6785.1665
$ws_CRP.Range("K122").Value = 6374.625
$ws_CRP.Range("L122").Value = 20355.4995
$ws_CRP.Range("M122").Value = -3924.625
$ws_CRP.Range("N122").Value = -25255.4995

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 573.0476
$ws_CUL.Range("I5").Value = 520.41174
$ws_CUL.Range("J5").Value = 796.75
$ws_CUL.Range("K5").Value = 1561.23522
$ws_CUL.Range("L5").Value = 2390.25
$ws_CUL.Range("M5").Value = -1449.23522
$ws_CUL.Range("N5").Value = -2614.25

$ws_CUL.Range("H81").Value = 2743.8572
$ws_CUL.Range("J81").Value = 3001.1667
$ws_CUL.Range("L81").Value = 9003.500100000001
$ws_CUL.Range("N81").Value = -11249.5001

$ws_CUL.Range("H84").Value = 2743.8572
$ws_CUL.Range("J84").Value = 3001.1667
$ws_CUL.Range("L84").Value = 27010.5003
$ws_CUL.Range("N84").Value = -38242.5003

$ws_CUL.Range("H122").Value = 843.44446
$ws_CUL.Range("I122").Value = 479.8
$ws_CUL.Range("J122").Value = 983.3077
$ws_CUL.Range("K122").Value = 4318.2
$ws_CUL.Range("L122").Value = 8849.7693
$ws_CUL.Range("M122").Value = -1868.2
$ws_CUL.Range("N122").Value = -13749.7693

$ws_CUL.Range("H132").Value = 10960.3
$ws_CUL.Range("J132").Value = 25874.75
$ws_CUL.Range("L132").Value = 232872.75
$ws_CUL.Range("N132").Value = -237932.75

$ws_CUL.Range("H135").Value = 573.0476
$ws_CUL.Range("I135").Value = 520.41174
$ws_CUL.Range("J135").Value = 796.75
$ws_CUL.Range("K135").Value = 4683.70566
$ws_CUL.Range("L135").Value = 7170.75
$ws_CUL.Range("M135").Value = -2148.70566
$ws_CUL.Range("N135").Value = -12240.75

$ws_CUL.Range("H137").Value = 3507.36
$ws_CUL.Range("I137").Value = 1755.5834
$ws_CUL.Range("J137").Value = 5124.385
$ws_CUL.Range("K137").Value = 5266.7502
$ws_CUL.Range("L137").Value = 15373.155
$ws_CUL.Range("M137").Value = -166.7502000000004
$ws_CUL.Range("N137").Value = -25573.155

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 1235.44
$ws_GSM.Range("I97").Value = 1166.95
$ws_GSM.Range("K97").Value = 1166.95
$ws_GSM.Range("M97").Value = -670.95

$ws_GSM.Range("H102").Value = 2270.8572
$ws_GSM.Range("I102").Value = 3224.25
$ws_GSM.Range("K102").Value = 3224.25
$ws_GSM.Range("M102").Value = -1602.25

$ws_GSM.Range("H122").Value = 1749.5
$ws_GSM.Range("I122").Value = 1685
$ws_GSM.Range("K122").Value = 5055
$ws_GSM.Range("M122").Value = -2605

$ws_GSM.Range("H132").Value = 1284408.8
$ws_GSM.Range("I132").Value = 1481201.8
$ws_GSM.Range("J132").Value = 5253.75
$ws_GSM.Range("K132").Value = 4443605.4
$ws_GSM.Range("L132").Value = 15761.25
$ws_GSM.Range("M132").Value = -4441075.4
$ws_GSM.Range("N132").Value = -20821.25

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 1017.0417
$ws_LTW.Range("I22").Value = 505.73685
$ws_LTW.Range("K22").Value = 505.73685
$ws_LTW.Range("M22").Value = -210.73685

$ws_LTW.Range("H27").Value = 1017.0417
$ws_LTW.Range("I27").Value = 505.73685
$ws_LTW.Range("K27").Value = 505.73685
$ws_LTW.Range("M27").Value = -398.73685

$ws_LTW.Range("H40").Value = 10093.818
$ws_LTW.Range("I40").Value = 10270.4
$ws_LTW.Range("J40").Value = 9715.429
$ws_LTW.Range("K40").Value = 10270.4
$ws_LTW.Range("L40").Value = 9715.429
$ws_LTW.Range("M40").Value = -10134.4
$ws_LTW.Range("N40").Value = -9987.429

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H70").Value = 46996.332
$ws_WVR.Range("J70").Value = 46996.332
$ws_WVR.Range("L70").Value = 46996.332
$ws_WVR.Range("N70").Value = -47626.332

$ws_WVR.Range("H73").Value = 46996.332
$ws_WVR.Range("J73").Value = 46996.332
$ws_WVR.Range("L73").Value = 46996.332
$ws_WVR.Range("N73").Value = -49180.332

$ws_WVR.Range("H96").Value = 5251.55
$ws_WVR.Range("I96").Value = 1700.7778
$ws_WVR.Range("J96").Value = 8156.727
$ws_WVR.Range("K96").Value = 1700.7778
$ws_WVR.Range("L96").Value = 8156.727
$ws_WVR.Range("M96").Value = -327.7778000000001
$ws_WVR.Range("N96").Value = -10902.727

$ws_WVR.Range("H132").Value = 1378.3948
$ws_WVR.Range("I132").Value = 1687.8096
$ws_WVR.Range("J132").Value = 996.17645
$ws_WVR.Range("K132").Value = 5063.4288
$ws_WVR.Range("L132").Value = 2988.52935
$ws_WVR.Range("M132").Value = -2533.4288
$ws_WVR.Range("N132").Value = -8048.529350000001

$ws_WVR.Range("H136").Value = 17922536
$ws_WVR.Range("I136").Value = 29240836
$ws_WVR.Range("J136").Value = 1892.9166
$ws_WVR.Range("K136").Value = 87722508
$ws_WVR.Range("L136").Value = 5678.7498
$ws_WVR.Range("M136").Value = -87719958
$ws_WVR.Range("N136").Value = -10778.7498
